$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write order follows the original author's entry sequence so the rebuilt
# sharedStrings table lands on the same new indices (81-85) as the target.

# ---- Row 35 / Row 36 column A: "Tar" (new, shared by both rows) ----
$ws.Range("A35").Value = "Tar"

# ---- Row 35 column C: tar cvf/tf/xvf cheat-sheet ----
$ws.Range("C35").Value = "`$ tar cvf wallpaper.tar WallPaper    //package all content in WallPaper folder`n`$ tar tf wallpaper.tar  //display content in pack`n`$ tar xvf wallpaper.tar   //extract pack"
$ws.Range("C35").WrapText = $true
$ws.Rows.Item(35).RowHeight = 63

# ---- Row 36 column A: "Tar" (reuses the same shared string) ----
$ws.Range("A36").Value = "Tar"

# ---- Row 36 column B: "Package option" ----
$ws.Range("B36").Value = "Package option"

# ---- Row 35 column B: "example" ----
$ws.Range("B35").Value = "example"

# ---- Row 36 column C: rich-text explanation of cvf / tf / xvf ----
$run1 = "create a package: cvf"
$run2 = "`nc–create create a new archive`nv–verbose verbosely list files processed`nf–file=ARCHIVE use archive file or device ARCHIVE`n"
$run3 = "display the content: tf"
$run4 = "`nt–list list the contents of an archive`nf–file=ARCHIVE use archive file or device ARCHIVE`n"
$run5 = "extract the package: xvf"
$run6 = "`nx–extract, –get extract files from an archive`nv–verbose verbosely list files processed`nf–file=ARCHIVE use archive file or device ARCHIVE"
$full36 = $run1 + $run2 + $run3 + $run4 + $run5 + $run6
$ws.Range("C36").Value = $full36
$ws.Range("C36").WrapText = $true

$pos = 1
$len1 = $run1.Length
$ws.Range("C36").Characters($pos, $len1).Font.Bold = $true
$pos = $pos + $len1
$len2 = $run2.Length
$c = $ws.Range("C36").Characters($pos, $len2)
$c.Font.Name = "Calibri"
$c.Font.Size = 12
$pos = $pos + $len2
$len3 = $run3.Length
$ws.Range("C36").Characters($pos, $len3).Font.Bold = $true
$pos = $pos + $len3
$len4 = $run4.Length
$c = $ws.Range("C36").Characters($pos, $len4)
$c.Font.Name = "Calibri"
$c.Font.Size = 12
$pos = $pos + $len4
$len5 = $run5.Length
$ws.Range("C36").Characters($pos, $len5).Font.Bold = $true
$pos = $pos + $len5
$len6 = $run6.Length
$c = $ws.Range("C36").Characters($pos, $len6)
$c.Font.Name = "Calibri"
$c.Font.Size = 12
$pos = $pos + $len6

$ws.Rows.Item(36).RowHeight = 173.25

# ---- view state: mirror scrolled position / selection from the authored edit ----
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 32
$win.ScrollColumn = 1
$ws.Range("B36").Select()

Write-Output "done"
